$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 15 middle detail rows (old rows 19-33) so only 4 employee rows remain
# (rows 16-18 stay in place, the old last row - 34 - shifts up to become row 19,
# and the footer rows - 39/40 - shift up to become 24/25).
$ws.Range("19:33").Delete()

# --- Header summary values ---
$ws.Range("E11").Value = 267480
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 2

# --- Updated detail rows (B stays "CC" on every row) ---
# Row 16: CINDY PATRICIA BURGOS PUELLO
$ws.Range("C16").Value = "1143354109"
$ws.Range("D16").Value = "CINDY PATRICIA BURGOS PUELLO"
$ws.Range("E16").Value = "1811"
$ws.Range("F16").Value = 5600
$ws.Range("G16").Value = 1400000

# Row 17: GABRIEL RAMON VELAZCO SALAZAR
$ws.Range("C17").Value = "73139841"
$ws.Range("D17").Value = "GABRIEL RAMON VELAZCO SALAZAR"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18: MOISES DE JESUS RUIZ PACHECO
$ws.Range("C18").Value = "7961217"
$ws.Range("D18").Value = "MOISES DE JESUS RUIZ PACHECO"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19: WILLIAM ALEJANDRO OLARTE RODRIGUEZ
$ws.Range("C19").Value = "11442202"
$ws.Range("D19").Value = "WILLIAM ALEJANDRO OLARTE RODRIGUEZ"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 148000
$ws.Range("G19").Value = 3700000
